$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / non-numeric-looking values
$ws.Range("D2").Value = "27.039.81"
$ws.Range("E2").Value = "  +2.43%  "
$ws.Range("D3").Value = "1.675.15"
$ws.Range("E3").Value = "  +3.51%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("E5").Value = "  +1.54%  "
$ws.Range("E6").Value = "  +2.17%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +2.96%  "
$ws.Range("E9").Value = "  +1.66%  "
$ws.Range("E10").Value = "  +5.47%  "
$ws.Range("E11").Value = "  +4.97%  "
$ws.Range("D12").Value = "1.911.71"
$ws.Range("E12").Value = "  +3.51%  "
$ws.Range("D13").Value = "1.675.76"
$ws.Range("E13").Value = "  +3.49%  "
$ws.Range("E14").Value = "  +1.87%  "
$ws.Range("E15").Value = "  +2.86%  "
$ws.Range("E16").Value = "  +3.20%  "
$ws.Range("D17").Value = "27.079.86"
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("D19").Value = "0.0₃0738"
$ws.Range("E19").Value = "  +1.80%  "
$ws.Range("E20").Value = "  -0.76%  "
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("E22").Value = "  +4.22%  "
$ws.Range("E23").Value = "  +2.24%  "
$ws.Range("E24").Value = "  +2.77%  "
$ws.Range("E25").Value = "  -0.90%  "
$ws.Range("E26").Value = "  +1.48%  "
$ws.Range("E27").Value = "  +0.59%  "
$ws.Range("E28").Value = "  +2.78%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("E30").Value = "  +0.55%  "
$ws.Range("E31").Value = "  +1.82%  "
$ws.Range("E32").Value = "  +2.36%  "
$ws.Range("D33").Value = "1.472.12"
$ws.Range("E34").Value = "  +5.44%  "
$ws.Range("E35").Value = "  +6.32%  "
$ws.Range("E36").Value = "  -0.46%  "
$ws.Range("E37").Value = "  +1.38%  "
$ws.Range("E38").Value = "  +8.18%  "
$ws.Range("E39").Value = "  +2.38%  "
$ws.Range("E40").Value = "  +3.68%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("E42").Value = "  +11.94%  "
$ws.Range("E43").Value = "  +3.52%  "
$ws.Range("E44").Value = "  +8.13%  "
$ws.Range("D45").Value = "1.822.07"
$ws.Range("E45").Value = "  +3.59%  "
$ws.Range("E46").Value = "  +2.07%  "
$ws.Range("E47").Value = "  -0.27%  "
$ws.Range("E48").Value = "  +2.35%  "
$ws.Range("E49").Value = "  +4.43%  "
$ws.Range("E50").Value = "  +1.25%  "
$ws.Range("E51").Value = "  +2.41%  "

# Numeric-looking values that must remain stored as text:
# force text format, assign, then clear the format so no residual styling remains
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.23"
$ws.Range("D5").ClearFormats()
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.22"
$ws.Range("D10").ClearFormats()
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.09"
$ws.Range("D14").ClearFormats()
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.90"
$ws.Range("D16").ClearFormats()
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "236.75"
$ws.Range("D18").ClearFormats()
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.23"
$ws.Range("D24").ClearFormats()
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.70"
$ws.Range("D25").ClearFormats()
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.61"
$ws.Range("D35").ClearFormats()
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.12"
$ws.Range("D40").ClearFormats()
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").ClearFormats()
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.02"
$ws.Range("D42").ClearFormats()
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "66.66"
$ws.Range("D44").ClearFormats()
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.778"
$ws.Range("D46").ClearFormats()
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.33"
$ws.Range("D47").ClearFormats()
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.66"
$ws.Range("D51").ClearFormats()
